$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refreshed crypto snapshot: map of cell address -> new value
$updates = [ordered]@{
    'D2' = '28.520.45'
    'D3' = '1.830.20'
    'E3' = '  +2.10%  '
    'D4' = '1.001'
    'E4' = '  +0.02%  '
    'D5' = '315.50'
    'E5' = '  -0.18%  '
    'E6' = '  -0.02%  '
    'D7' = '0.5072'
    'E7' = '  -5.31%  '
    'D8' = '0.3923'
    'E8' = '  +2.41%  '
    'D9' = '0.07724'
    'E9' = '  +3.96%  '
    'D10' = '41.92'
    'E10' = '  +1.12%  '
    'E11' = '  +2.70%  '
    'D12' = '21.06'
    'E12' = '  +3.71%  '
    'D13' = '6.261'
    'E13' = '  +1.01%  '
    'E14' = '  +0.03%  '
    'D15' = '7.546'
    'E15' = '  +1.55%  '
    'D16' = '1.820.33'
    'E16' = '  +1.47%  '
    'D17' = '93.53'
    'E17' = '  +5.88%  '
    'D18' = '0.00001083'
    'E18' = '  +2.37%  '
    'D19' = '0.06615'
    'E19' = '  +1.43%  '
    'D20' = '17.76'
    'E20' = '  +2.64%  '
    'E21' = '  +0.01%  '
    'D22' = '6.134'
    'E22' = '  +2.89%  '
    'D23' = '28.530.91'
    'E23' = '  +2.54%  '
    'D24' = '11.14'
    'E24' = '  +0.33%  '
    'E25' = '  +7.63%  '
    'D26' = '156.95'
    'E26' = '  +0.22%  '
    'D27' = '20.65'
    'E27' = '  +2.27%  '
    'B28' = 'LidoDAOToken'
    'C28' = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
    'D28' = '2.429'
    'E28' = '  +4.56%  '
    'B29' = 'WrappedliquidstakedEther2.0'
    'C29' = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
    'D29' = '2.035.83'
    'E29' = '  +1.77%  '
    'D30' = '125.20'
    'E30' = '  +3.15%  '
    'D31' = '1.134'
    'E31' = '  +2.34%  '
    'D32' = '0.1091'
    'E32' = '  -0.03%  '
    'D33' = '5.663'
    'E33' = '  +2.85%  '
    'D34' = '3.667'
    'E34' = '  +0.48%  '
    'D35' = '0.07106'
    'E35' = '  +1.79%  '
    'D36' = '0.2217'
    'E36' = '  +1.32%  '
    'D37' = '8.983'
    'E37' = '  +6.83%  '
    'D38' = '0.02323'
    'E38' = '  +2.17%  '
    'D39' = '5.129'
    'E39' = '  +1.66%  '
    'D40' = '0.6243'
    'E40' = '  +2.35%  '
    'D41' = '11.22'
    'E41' = '  -1.35%  '
    'D42' = '1.190'
    'E42' = '  +2.60%  '
    'E43' = '  -0.13%  '
    'D44' = '1.397'
    'E44' = '  -1.02%  '
    'D45' = '13.43'
    'E45' = '  +1.67%  '
    'D46' = '0.5900'
    'E46' = '  +3.57%  '
    'D47' = '3.720'
    'E47' = '  +1.09%  '
    'D48' = '124.33'
    'E48' = '  -0.67%  '
    'D49' = '1.972'
    'E49' = '  +3.37%  '
    'E50' = '  +1.20%  '
    'D51' = '0.06926'
    'E51' = '  +2.06%  '
}

foreach ($addr in $updates.Keys) {
    $newVal = $updates[$addr]
    $cell = $ws.Range($addr)

    # Plain decimal strings (e.g. "1.001", "21.06") would otherwise be auto-
    # converted by Excel into numbers, losing the original text formatting of
    # the source price column. Force those specific cells to Text, write the
    # value, then restore the original cell style so no formatting is added.
    if ($newVal -match '^\d+(\.\d+)?$') {
        $origStyle = $cell.Style
        $cell.NumberFormat = '@'
        $cell.Value = $newVal
        $cell.Style = $origStyle
    } else {
        $cell.Value = $newVal
    }
}
